$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.479.71"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "1.674.01"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("D5").Value = "'219.84"
$ws.Range("E5").Value = "  +1.63%  "
$ws.Range("D6").Value = "'0.5320"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("D8").Value = "'0.2694"
$ws.Range("E8").Value = "  +3.07%  "
$ws.Range("D9").Value = "'0.06392"
$ws.Range("E9").Value = "  +0.51%  "
$ws.Range("D10").Value = "'21.75"
$ws.Range("E10").Value = "  +4.32%  "
$ws.Range("D11").Value = "'0.07798"
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "'4.499"
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("D13").Value = "1.673.52"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").Value = "'0.5573"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "0.0₅8340"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "'65.66"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "26.502.15"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "'4.761"
$ws.Range("E19").Value = "  +0.44%  "
$ws.Range("D21").Value = "'10.33"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").Value = "'6.316"
$ws.Range("E22").Value = "  +1.34%  "
$ws.Range("D23").Value = "'1.003"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +5.31%  "
$ws.Range("D25").Value = "'139.00"
$ws.Range("E25").Value = "  -4.97%  "
$ws.Range("D26").Value = "'7.416"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("E27").Value = "  +2.77%  "
$ws.Range("D28").Value = "'1.427"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").Value = "'0.06270"
$ws.Range("E29").Value = "  +5.18%  "
$ws.Range("D30").Value = "'1.292"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").Value = "'3.606"
$ws.Range("E31").Value = "  +5.74%  "
$ws.Range("D32").Value = "'3.428"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").Value = "'1.689"
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("D34").Value = "'1.011"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").Value = "'0.6149"
$ws.Range("E35").Value = "  +9.03%  "
$ws.Range("D36").Value = "'2.423"
$ws.Range("E36").Value = "  +1.22%  "
$ws.Range("D37").Value = "'2.782"
$ws.Range("E37").Value = "  +1.04%  "
$ws.Range("E38").Value = "  +4.72%  "
$ws.Range("D39").Value = "'0.01618"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "1.094.50"
$ws.Range("D41").Value = "'0.8590"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "'100.62"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").Value = "1.820.71"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("E45").Value = "  +1.39%  "
$ws.Range("D46").Value = "'58.66"
$ws.Range("E46").Value = "  +4.91%  "
$ws.Range("D47").Value = "'8.176"
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = "'1.516"
$ws.Range("E49").Value = "  +9.36%  "
$ws.Range("D50").Value = "'0.05194"
$ws.Range("E50").Value = "  +0.33%  "
$ws.Range("D51").Value = "'6.010"
$ws.Range("E51").Value = "  +1.18%  "
